# Generate Report for Handoff
# Rotates the e2e test-doc GUID (496dd6b7-...) to a fresh one
# (973110e4-31c7-486c-8a6e-43bb20f44b29) and refreshes the handoff
# timestamps across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newGuid = "973110e4-31c7-486c-8a6e-43bb20f44b29"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-06 17:41:31"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Column -eq 2) {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.2ce3fd87bca1f15ba9c25f8cececda0f8f9c10c9.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-06 17:41:26"
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.Range.Column -eq 9) {
        $hl.Delete()
    } elseif ($hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

$wsZh.Columns.Item(9).ColumnWidth = 17.76
$wsZh.Columns.Item(10).ColumnWidth = 20.76

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.2ce3fd87bca1f15ba9c25f8cececda0f8f9c10c9.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-06 17:41:31"
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.Range.Column -eq 9) {
        $hl.Delete()
    } elseif ($hl.Range.Column -eq 1) {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

$wsDe.Columns.Item(9).ColumnWidth = 17.76
$wsDe.Columns.Item(10).ColumnWidth = 20.76
